# preparing eighth run: simple QFV
# - return Q3 to nominal interval
# - adjust V1 which maxed out in initial test

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Q3 interval (D10:G10 define "Q3 ~ uniform(F10, G10)"): return upper bound
# to its nominal value of 6.5 (was temporarily narrowed to 6 for a prior run).
$ws.Range("G10").Value = 6.5

# V1 interval (D15:G15): widen the upper bound from 0.1 to 0.2 - V1hyper was
# maxing out (saturating near 1) against the too-narrow original range.
$ws.Range("G15").Value = 0.2

# V2 interval (D16:G16): shift its lower bound up to match V1's new upper
# bound (0.1 -> 0.2) so the two intervals stay contiguous.
$ws.Range("F16").Value = 0.2

# Recalculate all of the dependent formula cells (H/U/V/W/X columns, the
# row-22 summary concatenations, etc.) against the updated inputs.
$excel.CalculateFullRebuild()

# View state: move off the scrolled-right K1 anchor back toward the left of
# the sheet and land the selection on A17 for the next run.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A17").Select()

# Best-effort restore of the saved window chrome height (application-level
# window sizing; harmless if the host does not persist it).
$excel.ActiveWindow.Height = 8700
$wb.Windows.Item(1).Height = 8700
